$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = '48.247.33'
$ws.Range("E2").Value = '  +0.26%  '
$ws.Range("D3").Value = '2.501.78'
$ws.Range("E3").Value = '  -0.85%  '
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").Value = '317.99'
$ws.Range("E5").Value = '  -1.72%  '
$ws.Range("D6").Value = '106.19'
$ws.Range("E6").Value = '  -2.46%  '
$ws.Range("E7").Value = '  -1.54%  '
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("D9").Value = '0.538'
$ws.Range("E9").Value = '  -3.28%  '
$ws.Range("D10").Value = '38.99'
$ws.Range("E10").Value = '  -3.85%  '
$ws.Range("D11").Value = '20.29'
$ws.Range("E11").Value = '  -0.07%  '
$ws.Range("E12").Value = '  -2.26%  '
$ws.Range("E13").Value = '  -0.05%  '
$ws.Range("D14").Value = '7.11'
$ws.Range("E14").Value = '  -2.11%  '
$ws.Range("E15").Value = '  -0.95%  '
$ws.Range("D16").Value = '2.497.37'
$ws.Range("E16").Value = '  -1.14%  '
$ws.Range("D17").Value = '0.831'
$ws.Range("E17").Value = '  -3.19%  '
$ws.Range("D18").Value = '48.079.68'
$ws.Range("E18").Value = '  +0.20%  '
$ws.Range("D19").Value = '3.01'
$ws.Range("E19").Value = '  +11.62%  '
$ws.Range("E20").Value = '  -3.09%  '
$ws.Range("D21").Value = '6.61'
$ws.Range("E21").Value = '  -0.47%  '
$ws.Range("D22").Value = '0.0₃0931'
$ws.Range("E22").Value = '  -2.11%  '
$ws.Range("D23").Value = '71.18'
$ws.Range("E23").Value = '  -1.61%  '
$ws.Range("D24").Value = '268.93'
$ws.Range("E24").Value = '  -0.22%  '
$ws.Range("E25").Value = '  -2.29%  '
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  +0.22%  '
$ws.Range("D27").Value = '25.86'
$ws.Range("E27").Value = '  -1.22%  '
$ws.Range("D28").Value = '2.23'
$ws.Range("E28").Value = '  +1.18%  '
$ws.Range("E29").Value = '  -3.98%  '
$ws.Range("D30").Value = '0.140'
$ws.Range("E30").Value = '  -3.91%  '
$ws.Range("D31").Value = '34.67'
$ws.Range("E31").Value = '  -2.77%  '
$ws.Range("D32").Value = '49.35'
$ws.Range("E32").Value = '  -0.83%  '
$ws.Range("E33").Value = '  +0.00%  '
$ws.Range("D34").Value = '19.17'
$ws.Range("E34").Value = '  -3.86%  '
$ws.Range("E35").Value = '  -1.99%  '
$ws.Range("D36").Value = '0.0775'
$ws.Range("E36").Value = '  -2.52%  '
$ws.Range("D37").Value = '1.95'
$ws.Range("E37").Value = '  -2.33%  '
$ws.Range("E38").Value = '  -3.00%  '
$ws.Range("E39").Value = '  -3.54%  '
$ws.Range("D40").Value = '123.29'
$ws.Range("E40").Value = '  +3.49%  '
$ws.Range("B41").Value = 'Stellar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D41").Value = '0.111'
$ws.Range("E41").Value = '  -1.42%  '
$ws.Range("B42").Value = 'EnergySwap'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D42").Value = '22.30'
$ws.Range("E42").Value = '  +0.17%  '
$ws.Range("D43").Value = '2.22'
$ws.Range("E43").Value = '  +0.93%  '
$ws.Range("D45").Value = '2.003.76'
$ws.Range("E45").Value = '  -0.35%  '
$ws.Range("E46").Value = '  +0.80%  '
$ws.Range("E47").Value = '  +1.76%  '
$ws.Range("E48").Value = '  -2.72%  '
$ws.Range("D49").Value = '8.95'
$ws.Range("E49").Value = '  -2.17%  '
$ws.Range("D50").Value = '5.19'
$ws.Range("E50").Value = '  -1.21%  '
$ws.Range("E51").Value = '  -1.18%  '
